$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 40, shifting existing rows 40-152 down to 41-153
$ws.Rows(40).Insert()

# Populate the new row 40 with the new weekly data entry
$ws.Cells.Item(40, 1).Value = 10
$ws.Cells.Item(40, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(40, 3).Value = "La Araucanía"
$ws.Cells.Item(40, 4).Value = 45260
$ws.Cells.Item(40, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(40, 5).Value = 9
$ws.Cells.Item(40, 6).Value = "Fruta"
$ws.Cells.Item(40, 7).Value = 100101
$ws.Cells.Item(40, 8).Value = "Berries"
$ws.Cells.Item(40, 9).Value = 100101001
$ws.Cells.Item(40, 10).Value = "Arándano (blue)"
$ws.Cells.Item(40, 11).Value = "Sin especificar"
$ws.Cells.Item(40, 12).Value = "Primera"
$ws.Cells.Item(40, 13).Value = 1000
$ws.Cells.Item(40, 14).Value = 3300
$ws.Cells.Item(40, 15).Value = 3400
$ws.Cells.Item(40, 16).Value = 3350
$ws.Cells.Item(40, 17).Value = "$/kilo"
$ws.Cells.Item(40, 18).Value = "Región del Maule"
$ws.Cells.Item(40, 19).Value = 3350
$ws.Cells.Item(40, 20).Value = 1
